$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 ("Directive" / " ENVIRONMENT " / 1) entirely, shifting
# row 3 ("VT Biodiversity Project" / "PROJECT" / 2) up to become row 2.
$ws.Rows.Item(2).Delete()
